$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, $cellRef, $newValue)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}


$ws = $wb.Worksheets.Item("Commitments")
Set-TextValue $ws "B2" "CMT - 000957"
Set-TextValue $ws "H2" "10/21/2022"
Set-TextValue $ws "B3" "CMT - 000958"
Set-TextValue $ws "H3" "10/22/2022"
Set-TextValue $ws "B4" "CMT - 000959"
Set-TextValue $ws "H4" "10/28/2022"
Set-TextValue $ws "B5" "CMT - 000960"
Set-TextValue $ws "H5" "10/29/2022"
Set-TextValue $ws "B6" "CMT - 000961"
Set-TextValue $ws "H6" "11/1/2022"
Set-TextValue $ws "B7" "CMT - 000962"
Set-TextValue $ws "H7" "11/4/2022"
Set-TextValue $ws "B8" "CMT - 000963"
Set-TextValue $ws "H8" "11/9/2022"
Set-TextValue $ws "B9" "CMT - 000957"
Set-TextValue $ws "H9" "10/21/2022"

$ws = $wb.Worksheets.Item("CapitalCall")
Set-TextValue $ws "B2" "CC-0627"
Set-TextValue $ws "C2" "DD-0217"
Set-TextValue $ws "D2" "CMT - 000961"
Set-TextValue $ws "I2" "10/25/2022"
Set-TextValue $ws "J2" "10/26/2022"
Set-TextValue $ws "K2" "2000000"
Set-TextValue $ws "L2" "10/17/2022"
Set-TextValue $ws "B3" "CC-0628"
Set-TextValue $ws "C3" "DD-0217"
Set-TextValue $ws "D3" "CMT - 000960"
Set-TextValue $ws "I3" "10/25/2022"
Set-TextValue $ws "J3" "10/26/2022"
Set-TextValue $ws "B4" "CC-0629"
Set-TextValue $ws "C4" "DD-0217"
Set-TextValue $ws "D4" "CMT - 000957"
Set-TextValue $ws "I4" "10/25/2022"
Set-TextValue $ws "J4" "10/26/2022"
Set-TextValue $ws "L4" "10/25/2022"
Set-TextValue $ws "B5" "CC-0630"
Set-TextValue $ws "C5" "DD-0217"
Set-TextValue $ws "D5" "CMT - 000958"
Set-TextValue $ws "I5" "10/25/2022"
Set-TextValue $ws "J5" "10/26/2022"
Set-TextValue $ws "B6" "CC-0631"
Set-TextValue $ws "C6" "DD-0217"
Set-TextValue $ws "D6" "CMT - 000959"
Set-TextValue $ws "I6" "10/25/2022"
Set-TextValue $ws "J6" "10/26/2022"

$ws = $wb.Worksheets.Item("FundDrawdown")
Set-TextValue $ws "B2" "DD-0217"
Set-TextValue $ws "E2" "10/25/2022"
Set-TextValue $ws "F2" "10/26/2022"

$ws = $wb.Worksheets.Item("FundDistribution")
Set-TextValue $ws "B2" "FD-0230"

$ws = $wb.Worksheets.Item("InvestorDistribution")
Set-TextValue $ws "B2" "FD-0230"
Set-TextValue $ws "C2" "ID-0633"
Set-TextValue $ws "D2" "CMT - 000961"
Set-TextValue $ws "B3" "FD-0230"
Set-TextValue $ws "C3" "ID-0634"
Set-TextValue $ws "D3" "CMT - 000960"
Set-TextValue $ws "B4" "FD-0230"
Set-TextValue $ws "C4" "ID-0635"
Set-TextValue $ws "D4" "CMT - 000957"
Set-TextValue $ws "B5" "FD-0230"
Set-TextValue $ws "C5" "ID-0636"
Set-TextValue $ws "D5" "CMT - 000958"
Set-TextValue $ws "B6" "FD-0230"
Set-TextValue $ws "C6" "ID-0637"
Set-TextValue $ws "D6" "CMT - 000959"

$ws = $wb.Worksheets.Item("CustomEmailFolder")
Set-TextValue $ws "B2" "PETestEmailFolder23854"
Set-TextValue $ws "C2" "PETestCustomEmailTemplate37534"

$ws = $wb.Worksheets.Item("Report")
Set-TextValue $ws "B2" "CustomReportFolder7404"
Set-TextValue $ws "C2" "CustomReport53088"

$ws = $wb.Worksheets.Item("Contacts")
Set-TextValue $ws "E2" "navatariptesting+59874@gmail.com"
Set-TextValue $ws "E3" "navatariptesting+45689@gmail.com"
Set-TextValue $ws "E4" "navatariptesting+41597@gmail.com"
Set-TextValue $ws "E6" "navatariptesting+25845@gmail.com"
Set-TextValue $ws "E7" "navatariptesting+27632@gmail.com"

$ws = $wb.Worksheets.Item("Funds")
Set-TextValue $ws "F2" "8.25E8"
Set-TextValue $ws "H2" "10/21/2022"
Set-TextValue $ws "H3" "10/21/2022"
Set-TextValue $ws "F4" "1.12E8"
Set-TextValue $ws "H4" "10/21/2022"
